# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before the "总计" sheet, and
#    populate it with the fund-holdings detail rows (same layout as the
#    other quarterly sheets: 2020-Q4 .. 2021-Q3).
# 2. Insert a new summary row at the top of the "总计" sheet's data
#    (right under the header) for the "2022-Q1" quarter, and keep the
#    running index column (A) renumbered 0..4.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: add the "2022-Q1" worksheet just before "总计"
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

# Copy the header + first data row formatting from an existing quarterly
# sheet (2021-Q3) so the new sheet matches the established look (bold
# centered header cells, bordered index cell, etc.)
$template = $wb.Worksheets.Item("2021-Q3")
$template.Range("B1:H2").Copy($q1.Range("B1:H2"))
$template.Range("A2").Copy($q1.Range("A2"))

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "'006478"
$q1.Range("C2").Value = "长盛多因子策略优选股票"
$q1.Range("D2").Value = "'0.51"
$q1.Range("E2").Value = "'84.41"
$q1.Range("F2").Value = "'4.37"
$q1.Range("G2").Value = "'0.0223"
$q1.Range("H2").Value = 6

# ---------------------------------------------------------------------
# Step 2: insert the "2022-Q1" row into the "总计" sheet
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# The freshly inserted row is blank; pick up the same styling the other
# data rows use by copying from the row directly below (which now holds
# what used to be row 2).
$total.Range("A3:D3").Copy($total.Range("A2:D2"))

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.02

# Renumber the running index column (A) for the rows that got pushed
# down so it stays a contiguous 0-based sequence.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
